$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.318.84"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "2.015.90"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'251.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.81%  "
$ws.Range("E6").Value = "  -2.96%  "
$ws.Range("D7").Value = "'62.88"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +10.87%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'59.42"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.16%  "
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("D13").Value = "'0.907"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("E14").Value = "  +4.86%  "
$ws.Range("D15").Value = "2.307.99"
$ws.Range("E15").Value = "  -1.17%  "
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").Value = "'19.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +14.37%  "
$ws.Range("D18").Value = "2.017.55"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").Value = "36.244.93"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").Value = "0.0₃0861"
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("D22").Value = "'5.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.56%  "
$ws.Range("D23").Value = "'234.41"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.09%  "
$ws.Range("E24").Value = "  +18.14%  "
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").Value = "'2.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.84%  "
$ws.Range("E27").Value = "  +3.37%  "
$ws.Range("D28").Value = "'162.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.62%  "
$ws.Range("D29").Value = "'19.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("D30").Value = "'0.116"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +33.89%  "
$ws.Range("D31").Value = "'0.121"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("E32").Value = "  +3.39%  "
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0609"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.02%  "
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "'4.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.43%  "
$ws.Range("E36").Value = "  +12.62%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").Value = "'5.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +16.88%  "
$ws.Range("E40").Value = "  +14.46%  "
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("E42").Value = "  +1.32%  "
$ws.Range("E43").Value = "  +0.96%  "
$ws.Range("E44").Value = "  +3.06%  "
$ws.Range("E45").Value = "  +4.86%  "
$ws.Range("E46").Value = "  +7.22%  "
$ws.Range("D47").Value = "'94.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.54%  "
$ws.Range("D48").Value = "1.432.79"
$ws.Range("E48").Value = "  +4.96%  "
$ws.Range("D49").Value = "'2.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +15.90%  "
$ws.Range("E50").Value = "  -0.56%  "
$ws.Range("D51").Value = "'47.46"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.87%  "
